$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append (rows 10-13)
$newRows = @(
    @(42613.758136574077, -2, 49, 49, 49, 48, 10010, 17746, 1867, 196, 196, 13, 12, "Bag"),
    @(42613.88585648148, -4, 48, 49, 48, 48, 11659, 20804, 2226, 212, 217, 15, 14, "Bag"),
    @(42614.884247685186, -14, 48, 48, 48, 88, 10673, 14293, 1541, 156, 157, 1, 8, "Bag"),
    @(42615.884733796294, -10, 51, 48, 51, 77, 8968, 15134, 1468, 170, 162, 2, 7, "Bag")
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A: copy the formatted date cell above so it keeps the same date style,
    # then overwrite with the new date serial value.
    $ws.Range("A9").Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $rowData[0]

    for ($c = 2; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }

    # Column N: shared string "Bag"
    $ws.Cells.Item($r, 14).Value = $rowData[13]
}
